# Generate Report for Handback
# - Mark "Status" as handed back (was "Ready for handoff").
# - Fill in "Latest Target File" (hyperlink to the source .md), "Latest Handback File"
#   (the generated .xlf) and the real "Latest Handback DateTime" for each language /
#   source-file row, on both the zh-cn and de-de sheets.
# - Column widths widen to fit the newly-populated columns.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: the per-language Status columns (E, F) just reflect the new
# status text.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Hyperlinks.Add(
    $zhcn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e6a00591bcc012338b316af9647957ffd76e100/e2e/1a28e42c-aec6-4c23-9c82-6e13822b1ecc.md",
    "",
    "",
    "1a28e42c-aec6-4c23-9c82-6e13822b1ecc.md"
) | Out-Null
$zhcn.Range("J2").Value = "1a28e42c-aec6-4c23-9c82-6e13822b1ecc.fb0fc659a1327715d76d5775b73371ab2cd6cfcb.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-29 16:32:48"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e6a00591bcc012338b316af9647957ffd76e100/e2e/3f0088a1-f74c-44b1-87a3-00892a92b168.md",
    "",
    "",
    "3f0088a1-f74c-44b1-87a3-00892a92b168.md"
) | Out-Null
$zhcn.Range("J3").Value = "3f0088a1-f74c-44b1-87a3-00892a92b168.48959684265158ae7a2772d5b0180073fa3367cf.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-29 16:32:48"

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Hyperlinks.Add(
    $dede.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e6a00591bcc012338b316af9647957ffd76e100/e2e/1a28e42c-aec6-4c23-9c82-6e13822b1ecc.md",
    "",
    "",
    "1a28e42c-aec6-4c23-9c82-6e13822b1ecc.md"
) | Out-Null
$dede.Range("J2").Value = "1a28e42c-aec6-4c23-9c82-6e13822b1ecc.fb0fc659a1327715d76d5775b73371ab2cd6cfcb.de-de.xlf"
$dede.Range("K2").Value = "2016-08-29 16:32:55"

$dede.Hyperlinks.Add(
    $dede.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e6a00591bcc012338b316af9647957ffd76e100/e2e/3f0088a1-f74c-44b1-87a3-00892a92b168.md",
    "",
    "",
    "3f0088a1-f74c-44b1-87a3-00892a92b168.md"
) | Out-Null
$dede.Range("J3").Value = "3f0088a1-f74c-44b1-87a3-00892a92b168.48959684265158ae7a2772d5b0180073fa3367cf.de-de.xlf"
$dede.Range("K3").Value = "2016-08-29 16:32:55"

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
